# Automatic update of files.
#
# The data rows (2-96) of the single worksheet get reshuffled into a new
# order (an external re-generation of the "Avverkningsanmälningar" list),
# and the "Förändrad" (column C) timestamp that is stamped on every row is
# advanced by one day (46065 -> 46066 serial, i.e. 2026-02-12 -> 2026-02-13).
#
# Strategy: snapshot every data row's values (columns A:R, via Value2 so
# numeric/date precision is preserved exactly) and formulas (columns S:Z,
# via Formula so the HYPERLINK() formulas survive) BEFORE writing anything
# back, bump column C by one day in the snapshot, then write each row back
# out at its new target position according to the row-permutation map
# below (target row -> source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 96

# target row -> source row (within the original sheet, rows 2..96)
$rowMap = @{ 2=3; 3=2; 4=4; 5=5; 6=6; 7=7; 8=8; 9=9; 10=10; 11=12; 12=13; 13=11; 14=15; 15=20; 16=45; 17=58; 18=90; 19=62; 20=96; 21=36; 22=33; 23=29; 24=17; 25=71; 26=34; 27=43; 28=19; 29=24; 30=73; 31=52; 32=74; 33=46; 34=49; 35=48; 36=75; 37=77; 38=51; 39=57; 40=53; 41=27; 42=65; 43=66; 44=82; 45=23; 46=67; 47=41; 48=68; 49=79; 50=44; 51=59; 52=81; 53=80; 54=26; 55=39; 56=40; 57=21; 58=69; 59=56; 60=87; 61=31; 62=78; 63=89; 64=47; 65=88; 66=50; 67=91; 68=76; 69=70; 70=42; 71=14; 72=61; 73=93; 74=64; 75=94; 76=92; 77=95; 78=25; 79=22; 80=32; 81=38; 82=63; 83=55; 84=16; 85=18; 86=54; 87=37; 88=60; 89=72; 90=30; 91=35; 92=28; 93=83; 94=84; 95=85; 96=86 }

# --- 1. Snapshot every source row (values A:R + formulas S:Z) ---------
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $valArr = $ws.Range("A${r}:R${r}").Value2
    $fmlArr = $ws.Range("S${r}:Z${r}").Formula

    # Column C (3rd column of the A:R block) is the "Förändrad" date -
    # every row advances by exactly one day.
    if ($null -ne $valArr[1,3]) {
        $valArr[1,3] = $valArr[1,3] + 1
    }

    $snapshot[$r] = @{ val = $valArr; fml = $fmlArr }
}

# --- 2. Write every row back out at its new (target) position ---------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $rowMap[$r]
    $row = $snapshot[$src]
    $ws.Range("A${r}:R${r}").Value2 = $row.val
    $ws.Range("S${r}:Z${r}").Formula = $row.fml
}
